$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) - updates columns B,C,D,E,G for rows 2-31
$data = @{
    2 = @(0.6545652718822623, 0.002658071450198252, 0.7210945179870265, 0.5333859586016987, 1.911703819921186)
    3 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    4 = @(0.0006075818656279264, 0.002658071450198252, 3.223369029078222, 0.5333859586016987, 3.760020640995746)
    5 = @(1.445647641019636, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 3.005019366241741)
    6 = @(0.2881169905109251, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1.012145535086602)
    7 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    8 = @(0.1169995834814548, 0.002658071450198252, 3.223369029078222, 0.5333859586016987, 3.876412642611573)
    9 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.178645819794754)
    10 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    11 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    12 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    13 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 18.91276827552123)
    14 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
    15 = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 22.32281868886277)
    16 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    17 = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.412515779045154)
    18 = @(0.2881169905109251, 0.3048912486333797, 18.71679738969934, 13.86384647080068, 33.17365209964433)
    19 = @(0.1169995834814548, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.67637130870356)
    20 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    21 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    22 = @(1.445647641019636, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 2.433531715253719)
    23 = @(0.2881169905109251, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.169585166641744)
    24 = @(0.1169995834814548, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 2.998467759612273)
    25 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    26 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    27 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    28 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    29 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    30 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
    31 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B: TB
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C: d2S
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D: K
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E: IP
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G: sum
}
